$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44211
$ws.Cells.Item(2, 10).Value = 70
$ws.Cells.Item(2, 11).Value = 22000
$ws.Cells.Item(2, 12).Value = 24000
$ws.Cells.Item(2, 13).Value = 23143
$ws.Cells.Item(2, 15).Value = "Región del Maule"
$ws.Cells.Item(2, 16).Value = 926

# Row 3
$ws.Cells.Item(3, 4).Value = 44208
$ws.Cells.Item(3, 10).Value = 65
$ws.Cells.Item(3, 12).Value = 25000
$ws.Cells.Item(3, 13).Value = 23385
$ws.Cells.Item(3, 16).Value = 935

# Row 4
$ws.Cells.Item(4, 4).Value = 44193
$ws.Cells.Item(4, 10).Value = 30
$ws.Cells.Item(4, 11).Value = 35000
$ws.Cells.Item(4, 12).Value = 36000
$ws.Cells.Item(4, 13).Value = 35500
$ws.Cells.Item(4, 16).Value = 1420

# Row 5
$ws.Cells.Item(5, 4).Value = 44160
$ws.Cells.Item(5, 10).Value = 30
$ws.Cells.Item(5, 11).Value = 30000
$ws.Cells.Item(5, 12).Value = 30000
$ws.Cells.Item(5, 13).Value = 30000
$ws.Cells.Item(5, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(5, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(5, 16).Value = 1200

# Row 6
$ws.Cells.Item(6, 4).Value = 44160
$ws.Cells.Item(6, 8).Value = "Magnum"
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(6, 11).Value = 28000
$ws.Cells.Item(6, 12).Value = 28000
$ws.Cells.Item(6, 13).Value = 28000
$ws.Cells.Item(6, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(6, 16).Value = 1120

# Row 7
$ws.Cells.Item(7, 4).Value = 44201
$ws.Cells.Item(7, 10).Value = 33
$ws.Cells.Item(7, 11).Value = 26000
$ws.Cells.Item(7, 12).Value = 28000
$ws.Cells.Item(7, 13).Value = 27091
$ws.Cells.Item(7, 15).Value = "Región del Maule"
$ws.Cells.Item(7, 16).Value = 1084

# Row 8
$ws.Cells.Item(8, 4).Value = 44232
$ws.Cells.Item(8, 10).Value = 30
$ws.Cells.Item(8, 11).Value = 24000
$ws.Cells.Item(8, 12).Value = 25000
$ws.Cells.Item(8, 13).Value = 24500
$ws.Cells.Item(8, 16).Value = 980

# Row 9
$ws.Cells.Item(9, 4).Value = 44166
$ws.Cells.Item(9, 8).Value = "Magnum"
$ws.Cells.Item(9, 10).Value = 38
$ws.Cells.Item(9, 13).Value = 24526
$ws.Cells.Item(9, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(9, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(9, 16).Value = 981

# Row 10
$ws.Cells.Item(10, 4).Value = 44323
$ws.Cells.Item(10, 10).Value = 30
$ws.Cells.Item(10, 11).Value = 29000
$ws.Cells.Item(10, 12).Value = 30000
$ws.Cells.Item(10, 13).Value = 29500
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 1180

# Row 11
$ws.Cells.Item(11, 4).Value = 44281
$ws.Cells.Item(11, 10).Value = 30
$ws.Cells.Item(11, 11).Value = 31000
$ws.Cells.Item(11, 12).Value = 32000
$ws.Cells.Item(11, 13).Value = 31500
$ws.Cells.Item(11, 15).Value = "Región del Maule"
$ws.Cells.Item(11, 16).Value = 1260

# Row 12
$ws.Cells.Item(12, 4).Value = 44181
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 10).Value = 60
$ws.Cells.Item(12, 11).Value = 19500
$ws.Cells.Item(12, 12).Value = 20000
$ws.Cells.Item(12, 13).Value = 19750
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Región del Maule"
$ws.Cells.Item(12, 16).Value = 790

# Row 13
$ws.Cells.Item(13, 4).Value = 44246
$ws.Cells.Item(13, 10).Value = 60
$ws.Cells.Item(13, 11).Value = 24000
$ws.Cells.Item(13, 12).Value = 25000
$ws.Cells.Item(13, 13).Value = 24500
$ws.Cells.Item(13, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(13, 16).Value = 980

# Row 14
$ws.Cells.Item(14, 4).Value = 44203
$ws.Cells.Item(14, 10).Value = 50
$ws.Cells.Item(14, 11).Value = 22000
$ws.Cells.Item(14, 12).Value = 24000
$ws.Cells.Item(14, 13).Value = 23200
$ws.Cells.Item(14, 16).Value = 928

# Row 16
$ws.Cells.Item(16, 4).Value = 44174
$ws.Cells.Item(16, 11).Value = 21000
$ws.Cells.Item(16, 12).Value = 22000
$ws.Cells.Item(16, 13).Value = 21500
$ws.Cells.Item(16, 16).Value = 860

# Row 17
$ws.Cells.Item(17, 4).Value = 44236
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 22000
$ws.Cells.Item(17, 12).Value = 23000
$ws.Cells.Item(17, 13).Value = 22500
$ws.Cells.Item(17, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(17, 16).Value = 900

# Row 18
$ws.Cells.Item(18, 4).Value = 44159
$ws.Cells.Item(18, 8).Value = "Magnum"
$ws.Cells.Item(18, 10).Value = 47
$ws.Cells.Item(18, 11).Value = 27000
$ws.Cells.Item(18, 12).Value = 28000
$ws.Cells.Item(18, 13).Value = 27532
$ws.Cells.Item(18, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(18, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(18, 16).Value = 1101

# Row 19
$ws.Cells.Item(19, 4).Value = 44195
$ws.Cells.Item(19, 10).Value = 160
$ws.Cells.Item(19, 11).Value = 32000
$ws.Cells.Item(19, 12).Value = 33000
$ws.Cells.Item(19, 13).Value = 32500
$ws.Cells.Item(19, 16).Value = 1300

# Row 20
$ws.Cells.Item(20, 4).Value = 44209
$ws.Cells.Item(20, 10).Value = 90
$ws.Cells.Item(20, 11).Value = 23000
$ws.Cells.Item(20, 12).Value = 25000
$ws.Cells.Item(20, 13).Value = 23889
$ws.Cells.Item(20, 15).Value = "Región del Maule"
$ws.Cells.Item(20, 16).Value = 956

# Row 21
$ws.Cells.Item(21, 4).Value = 44250
$ws.Cells.Item(21, 10).Value = 120
$ws.Cells.Item(21, 11).Value = 22000
$ws.Cells.Item(21, 12).Value = 23000
$ws.Cells.Item(21, 13).Value = 22500
$ws.Cells.Item(21, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(21, 16).Value = 900

# Row 22
$ws.Cells.Item(22, 4).Value = 44334
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 10).Value = 30
$ws.Cells.Item(22, 11).Value = 31000
$ws.Cells.Item(22, 12).Value = 32000
$ws.Cells.Item(22, 13).Value = 31500
$ws.Cells.Item(22, 15).Value = "Región del Maule"
$ws.Cells.Item(22, 16).Value = 1260

# Row 23
$ws.Cells.Item(23, 4).Value = 44186
$ws.Cells.Item(23, 10).Value = 60
$ws.Cells.Item(23, 11).Value = 19000
$ws.Cells.Item(23, 12).Value = 20000
$ws.Cells.Item(23, 13).Value = 19500
$ws.Cells.Item(23, 16).Value = 780

# Row 24
$ws.Cells.Item(24, 4).Value = 44252
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 10).Value = 80
$ws.Cells.Item(24, 11).Value = 22000
$ws.Cells.Item(24, 12).Value = 23000
$ws.Cells.Item(24, 13).Value = 22500
$ws.Cells.Item(24, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(24, 16).Value = 900

# Row 25
$ws.Cells.Item(25, 4).Value = 44267
$ws.Cells.Item(25, 10).Value = 45
$ws.Cells.Item(25, 11).Value = 24000
$ws.Cells.Item(25, 12).Value = 25000
$ws.Cells.Item(25, 13).Value = 24333
$ws.Cells.Item(25, 16).Value = 973

# Row 26
$ws.Cells.Item(26, 4).Value = 44259
$ws.Cells.Item(26, 10).Value = 65
$ws.Cells.Item(26, 11).Value = 24000
$ws.Cells.Item(26, 12).Value = 25000
$ws.Cells.Item(26, 13).Value = 24538
$ws.Cells.Item(26, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(26, 16).Value = 982

# Row 27
$ws.Cells.Item(27, 4).Value = 44179
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 11).Value = 16000
$ws.Cells.Item(27, 12).Value = 17000
$ws.Cells.Item(27, 13).Value = 16500
$ws.Cells.Item(27, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región del Maule"
$ws.Cells.Item(27, 16).Value = 660

# Row 28
$ws.Cells.Item(28, 4).Value = 44249
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 10).Value = 60
$ws.Cells.Item(28, 11).Value = 21000
$ws.Cells.Item(28, 12).Value = 22000
$ws.Cells.Item(28, 13).Value = 21500
$ws.Cells.Item(28, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(28, 16).Value = 860

# Row 29
$ws.Cells.Item(29, 4).Value = 44168
$ws.Cells.Item(29, 11).Value = 15000
$ws.Cells.Item(29, 12).Value = 16000
$ws.Cells.Item(29, 13).Value = 15500
$ws.Cells.Item(29, 16).Value = 620

# Row 30
$ws.Cells.Item(30, 4).Value = 44272
$ws.Cells.Item(30, 10).Value = 42
$ws.Cells.Item(30, 13).Value = 22857
$ws.Cells.Item(30, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(30, 16).Value = 914

# Row 31
$ws.Cells.Item(31, 4).Value = 44218
$ws.Cells.Item(31, 10).Value = 80
$ws.Cells.Item(31, 11).Value = 24000
$ws.Cells.Item(31, 12).Value = 25000
$ws.Cells.Item(31, 13).Value = 24562
$ws.Cells.Item(31, 15).Value = "Región del Maule"
$ws.Cells.Item(31, 16).Value = 982

# Row 32
$ws.Cells.Item(32, 4).Value = 44210
$ws.Cells.Item(32, 10).Value = 70
$ws.Cells.Item(32, 11).Value = 23000
$ws.Cells.Item(32, 13).Value = 23857
$ws.Cells.Item(32, 15).Value = "Región del Maule"
$ws.Cells.Item(32, 16).Value = 954

# Row 34
$ws.Cells.Item(34, 4).Value = 44176
$ws.Cells.Item(34, 10).Value = 30
$ws.Cells.Item(34, 11).Value = 19000
$ws.Cells.Item(34, 12).Value = 20000
$ws.Cells.Item(34, 13).Value = 19500
$ws.Cells.Item(34, 16).Value = 780

# Row 35
$ws.Cells.Item(35, 4).Value = 44273
$ws.Cells.Item(35, 10).Value = 33
$ws.Cells.Item(35, 11).Value = 23000
$ws.Cells.Item(35, 12).Value = 24000
$ws.Cells.Item(35, 13).Value = 23455
$ws.Cells.Item(35, 16).Value = 938

# Row 36
$ws.Cells.Item(36, 4).Value = 44302
$ws.Cells.Item(36, 8).Value = "Magnum"
$ws.Cells.Item(36, 10).Value = 60
$ws.Cells.Item(36, 11).Value = 25000
$ws.Cells.Item(36, 12).Value = 26000
$ws.Cells.Item(36, 13).Value = 25500
$ws.Cells.Item(36, 16).Value = 1020

# Row 37
$ws.Cells.Item(37, 4).Value = 44161
$ws.Cells.Item(37, 8).Value = "Magnum"
$ws.Cells.Item(37, 10).Value = 47
$ws.Cells.Item(37, 11).Value = 28000
$ws.Cells.Item(37, 12).Value = 29000
$ws.Cells.Item(37, 13).Value = 28532
$ws.Cells.Item(37, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(37, 16).Value = 1141
